$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.348.56'
$ws.Range("E2").Value = '  -0.83%  '

$ws.Range("D3").Value = '2.773.51'
$ws.Range("E3").Value = '  +0.22%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '352.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.70%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.46%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.548'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.43%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("E9").Value = '  -0.64%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.56'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.84%  '

$ws.Range("E11").Value = '  +2.86%  '

$ws.Range("B12").Value = 'Dogecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0834'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.68%  '

$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.04'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.55%  '

$ws.Range("E14").Value = '  -0.49%  '

$ws.Range("D15").Value = '3.210.13'
$ws.Range("E15").Value = '  +0.26%  '

$ws.Range("D16").Value = '2.771.51'
$ws.Range("E16").Value = '  -0.70%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.918'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.14%  '

$ws.Range("D18").Value = '51.366.18'
$ws.Range("E18").Value = '  -0.72%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.58'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.40%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.09'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.84%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.08%  '

$ws.Range("D22").Value = '0.0₃0962'
$ws.Range("E22").Value = '  -1.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.81'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.00%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '265.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.09%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.07%  '

$ws.Range("E27").Value = '  -1.89%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.162'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +11.80%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.21'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.87%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.56'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.93%  '

$ws.Range("E31").Value = '  -0.89%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.15'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +8.34%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '51.94'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.64%  '

$ws.Range("E34").Value = '  -2.16%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.53'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.47%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0827'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.70%  '

$ws.Range("E37").Value = '  +0.03%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.32'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.78%  '

$ws.Range("E39").Value = '  -2.42%  '

$ws.Range("E40").Value = '  -1.46%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.52'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.69%  '

$ws.Range("E42").Value = '  -0.43%  '

$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.10'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.40%  '

$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '119.84'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.04%  '

$ws.Range("E45").Value = '  -2.30%  '

$ws.Range("D46").Value = '2.101.63'
$ws.Range("E46").Value = '  +2.07%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.25'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.68%  '

$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.31'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.79%  '

$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.40'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.89%  '

$ws.Range("B50").Value = 'SEI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.899'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.55%  '

$ws.Range("E51").Value = '  +9.00%  '
